$d = $word.ActiveDocument

# Locate the paragraph that contains the instructional "Summarize the overall
# results..." placeholder text under the "Summary:" heading.
$needle = "Summarize the overall results of the deep learning model"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$needle*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $target = $d.Paragraphs.Item($targetIndex)

    # Delete the whole paragraph, including its end-of-paragraph mark, which
    # removes the placeholder sentence entirely and merges the following
    # (empty) paragraph up into its slot.
    $target.Range.Delete()

    # The paragraph that used to follow the deleted one (also blue/italic,
    # but empty) now sits at $targetIndex. Strip the italic formatting from
    # its paragraph mark while keeping the blue color.
    $survivor = $d.Paragraphs.Item($targetIndex)
    $survivor.Range.Font.Italic = $false
    $survivor.Range.Font.ItalicBi = $false
}
